# Commit: "add residential/unclassified road filter file"
#
# Adds unit rehabilitation costs for the "residential" and "unclassified"
# highway classes to the "road" sheet, mirroring the existing
# secondary/tertiary paved+unpaved cost rows (same $/km figures), and
# leaves the workbook focused on the "road" tab (matching the author's
# final view state: activeTab moved from "rail" to "road").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("road")

# --- residential: paved / unpaved -----------------------------------
$ws.Range("A14").Value = "residential"
$ws.Range("B14").Value = "paved"
$ws.Range("C14").Value = 164059.77207298
$ws.Range("C14").NumberFormat = "#,##0.00"

$ws.Range("A15").Value = "residential"
$ws.Range("B15").Value = "unpaved"
$ws.Range("C15").Value = 23545.1019998726
$ws.Range("C15").NumberFormat = "#,##0.00"

# --- unclassified: paved / unpaved -----------------------------------
$ws.Range("A16").Value = "unclassified"
$ws.Range("B16").Value = "paved"
$ws.Range("C16").Value = 164059.77207298
$ws.Range("C16").NumberFormat = "#,##0.00"

$ws.Range("A17").Value = "unclassified"
$ws.Range("B17").Value = "unpaved"
$ws.Range("C17").Value = 23545.1019998726
$ws.Range("C17").NumberFormat = "#,##0.00"

# --- leave the "road" sheet focused/active, as in the final workbook -
$ws.Activate()
[void]$ws.Range("E22").Select()
